# Shift the horizontal offset of the explicitly-positioned
# "Content Placeholder 2" shape on slides 2-9 by -1 EMU
# (782166 EMU -> 782165 EMU; the y-offset, 1769591 EMU, stays the same).
#
# Notes on the conversion:
# - The PowerPoint COM object model reports/accepts Shape.Left/.Top in
#   points (1 pt = 12700 EMU) through a single-precision (float32) value,
#   and the underlying engine truncates (floors) when converting the
#   point value back to EMU on save. A literal "782165/12700" therefore
#   round-trips to 782164 EMU, one short of the intended 782165. The
#   literal 61.58783 sits safely inside the float32 bucket that truncates
#   to exactly 782165 EMU, so it is used instead of the "exact" quotient.
#
# Some slides (1, 3, 4) contain a *second* shape also named
# "Content Placeholder 2" that has no explicit <a:xfrm> of its own (it
# inherits its position from the slide layout, reporting Left=66 /
# Top=143.75 through COM). That inherited-position shape must be left
# untouched - only the shape that already carries the explicit
# ~61.588 pt / ~139.338 pt position (i.e. the one this diff targets)
# should be edited. We detect it by checking how close Top is to the
# explicit placeholder's known position.
$targetLeftPoints = 61.58783

$p = $ppt.ActivePresentation

for ($i = 2; $i -le 9; $i++) {
    $s = $p.Slides.Item($i)
    foreach ($sh in $s.Shapes) {
        if ($sh.Name -eq "Content Placeholder 2" -and [Math]::Abs($sh.Top - 139.3379) -lt 0.5) {
            $sh.Left = $targetLeftPoints
        }
    }
}
